# "clean up full_obs for use with modules dir"
#
# - sheet "1" (the last sheet in the workbook) gains a 3rd data row
# - a brand-new sheet "2" is appended after sheet "1", carrying the row
#   that used to live alone on sheet "1" (row 2), but evolved further
# - the new sheet "2" becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- 1. Append a new row to the existing "1" sheet ---------------------
$sheet1 = $wb.Worksheets.Item("1")
$sheet1.Range("A3").Value = "add(add(add(y, x), conditional(y, y)), conditional(conditional(y, x), conditional(y, y)))"
$sheet1.Range("B3").Value = -1905

# --- 2. Create the new "2" sheet right after "1" ------------------------
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "2"

$sheet2.Range("A1").Value = "ind"
$sheet2.Range("B1").Value = "fitness"
$sheet2.Range("A2").Value = "conditional(conditional(conditional(x, vel), add(y, vel)), conditional(vel, vel))"
$sheet2.Range("B2").Value = -854

# --- 3. Make sure the new sheet ends up as the active tab ---------------
$sheet2.Activate()
